$wb = $excel.ActiveWorkbook

# Leave behind a selection on the first sheet (matches prior user navigation)
$firstSheet = $wb.Worksheets.Item("Tir_235_50R24")
$firstSheet.Activate()
$firstSheet.Range("C27").Select()

# Copy the last sheet (Tir_145_70R13) to create a new sheet at the end, following
# the same layout convention used for the other tires
$srcSheet = $wb.Worksheets.Item("Tir_145_70R13")
$srcSheet.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Tir_430_50R38"

# Update the tire name / reference cells on the new sheet
$newSheet.Range("H3").Value = "Delft_430_50R38"
$newSheet.Range("H5").Value = "which('Truck_430_50R38.tir')"

# Activate the new sheet
$newSheet.Activate()
$newSheet.Range("H6").Select()
